$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 34 (existing rows 34-68 shift down to 35-69).
$ws.Rows(34).Insert()

# Populate the new row 34 with the new weekly record (same market/region/
# product metadata as the rest of the sheet, new date + price figures).
$ws.Range("A34").Value = 10
$ws.Range("B34").Value = "Vega Modelo de Temuco"
$ws.Range("C34").Value = "La Araucanía"
$ws.Range("D34").Value = 44966
$ws.Range("E34").Value = 9
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100107
$ws.Range("H34").Value = "Otros"
$ws.Range("I34").Value = 100107011
$ws.Range("J34").Value = "Tuna"
$ws.Range("K34").Value = "Sin especificar"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 150
$ws.Range("N34").Value = 22000
$ws.Range("O34").Value = 22000
$ws.Range("P34").Value = 22000
$ws.Range("Q34").Value = "$/caja 18 kilos"
$ws.Range("R34").Value = "Provincia de Los Andes"
$ws.Range("S34").Value = 1222
$ws.Range("T34").Value = 18
